$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 held a leftover/placeholder patient record; replace it with the
# corrected data (codigo, nombre, apellidos, edad). Sex (E13) and
# vacunado (F13) are left untouched.
#
# Numeric-looking identifiers (codigo / edad) are prefixed with a quote so
# Excel stores them as text (matching the rest of the sheet, where these
# columns are text-typed), then ClearFormats() drops the quote-prefix
# formatting Excel applies for that so the cell keeps the workbook's
# default style.
$ws.Range("A13").Value = "'23495867"
$ws.Range("A13").ClearFormats()

$ws.Range("B13").Value = "BENIGNA "

$ws.Range("C13").Value = "ARMAS  JIMENEZ"

$ws.Range("D13").Value = "'56"
$ws.Range("D13").ClearFormats()
